$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 10000
$ws.Range("J32").Value = 10000
$ws.Range("L32").Value = 10000
$ws.Range("N32").Value = -10652
$ws.Range("H86").Value = 3565.875
$ws.Range("I86").Value = 4296.0454
$ws.Range("J86").Value = 1959.5
$ws.Range("K86").Value = 4296.0454
$ws.Range("L86").Value = 1959.5
$ws.Range("M86").Value = -3173.0454
$ws.Range("N86").Value = -4205.5
$ws.Range("H89").Value = 3565.875
$ws.Range("I89").Value = 4296.0454
$ws.Range("J89").Value = 1959.5
$ws.Range("K89").Value = 21480.227
$ws.Range("L89").Value = 9797.5
$ws.Range("M89").Value = -15864.227
$ws.Range("N89").Value = -21029.5
$ws.Range("H116").Value = 4867.364
$ws.Range("I116").Value = 4016.1667
$ws.Range("K116").Value = 4016.1667
$ws.Range("M116").Value = -574.1667000000002
$ws.Range("H135").Value = 5144.591
$ws.Range("I135").Value = 5294.3335
$ws.Range("K135").Value = 47649.0015
$ws.Range("M135").Value = -45114.0015
$ws.Range("H137").Value = 11632.645
$ws.Range("I137").Value = 2541.7693
$ws.Range("K137").Value = 7625.3079
$ws.Range("M137").Value = -5075.3079
$ws.Range("H138").Value = 4248.1777
$ws.Range("I138").Value = 6696.5
$ws.Range("J138").Value = 3548.6572
$ws.Range("K138").Value = 20089.5
$ws.Range("L138").Value = 10645.9716
$ws.Range("M138").Value = -14949.5
$ws.Range("N138").Value = -20925.9716

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1513.0526
$ws.Range("I45").Value = 1493.2667
$ws.Range("J45").Value = 1587.25
$ws.Range("K45").Value = 1493.2667
$ws.Range("L45").Value = 1587.25
$ws.Range("M45").Value = -1116.2667
$ws.Range("N45").Value = -2341.25
$ws.Range("H61").Value = 520600.8
$ws.Range("I61").Value = 2847
$ws.Range("J61").Value = 3281954.5
$ws.Range("K61").Value = 2847
$ws.Range("L61").Value = 3281954.5
$ws.Range("M61").Value = -2635
$ws.Range("N61").Value = -3282378.5
$ws.Range("H74").Value = 6584.9556
$ws.Range("I74").Value = 2148.3171
$ws.Range("J74").Value = 52060.5
$ws.Range("K74").Value = 2148.3171
$ws.Range("L74").Value = 52060.5
$ws.Range("M74").Value = -1274.3171
$ws.Range("N74").Value = -53808.5
$ws.Range("H77").Value = 6584.9556
$ws.Range("I77").Value = 2148.3171
$ws.Range("J77").Value = 52060.5
$ws.Range("K77").Value = 10741.5855
$ws.Range("L77").Value = 260302.5
$ws.Range("M77").Value = -6373.585500000001
$ws.Range("N77").Value = -269038.5
$ws.Range("H110").Value = 3554.353
$ws.Range("I110").Value = 4494.4614
$ws.Range("J110").Value = 499
$ws.Range("K110").Value = 4494.4614
$ws.Range("L110").Value = 499
$ws.Range("M110").Value = -2449.4614
$ws.Range("N110").Value = -4589
$ws.Range("H122").Value = 770537.9
$ws.Range("I122").Value = 1001201.4
$ws.Range("J122").Value = 1659.6666
$ws.Range("K122").Value = 3003604.2
$ws.Range("L122").Value = 4978.9998
$ws.Range("M122").Value = -3001154.2
$ws.Range("N122").Value = -9878.9998
$ws.Range("H132").Value = 515054.5
$ws.Range("I132").Value = 2846.0588
$ws.Range("J132").Value = 2691940.5
$ws.Range("K132").Value = 8538.1764
$ws.Range("L132").Value = 8075821.5
$ws.Range("M132").Value = -6008.1764
$ws.Range("N132").Value = -8080881.5
$ws.Range("H135").Value = 212714.5
$ws.Range("J135").Value = 212714.5
$ws.Range("L135").Value = 212714.5
$ws.Range("N135").Value = -222854.5
$ws.Range("H136").Value = 520600.8
$ws.Range("I136").Value = 2847
$ws.Range("J136").Value = 3281954.5
$ws.Range("K136").Value = 8541
$ws.Range("L136").Value = 9845863.5
$ws.Range("M136").Value = -5991
$ws.Range("N136").Value = -9850963.5
$ws.Range("H139").Value = 82619.164
$ws.Range("J139").Value = 82619.164
$ws.Range("L139").Value = 82619.164
$ws.Range("N139").Value = -92899.164

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H113").Value = 120000
$ws.Range("I113").Value = 120000
$ws.Range("K113").Value = 120000
$ws.Range("M113").Value = -117830
$ws.Range("H134").Value = 11179.903
$ws.Range("I134").Value = 7799.6665
$ws.Range("J134").Value = 19030.773
$ws.Range("K134").Value = 23398.9995
$ws.Range("L134").Value = 57092.319
$ws.Range("M134").Value = -20863.9995
$ws.Range("N134").Value = -62162.319

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 6068.5264
$ws.Range("I16").Value = 6884.5
$ws.Range("J16").Value = 1716.6666
$ws.Range("K16").Value = 6884.5
$ws.Range("L16").Value = 1716.6666
$ws.Range("M16").Value = -6597.5
$ws.Range("N16").Value = -2290.6666
$ws.Range("H31").Value = 204036.5
$ws.Range("I31").Value = 403222.2
$ws.Range("J31").Value = 61761
$ws.Range("K31").Value = 403222.2
$ws.Range("L31").Value = 61761
$ws.Range("M31").Value = -402927.2
$ws.Range("N31").Value = -62351
$ws.Range("H34").Value = 204036.5
$ws.Range("I34").Value = 403222.2
$ws.Range("J34").Value = 61761
$ws.Range("K34").Value = 403222.2
$ws.Range("L34").Value = 61761
$ws.Range("M34").Value = -403020.2
$ws.Range("N34").Value = -62165
$ws.Range("H94").Value = 13791.75
$ws.Range("I94").Value = 13441.667
$ws.Range("J94").Value = 13908.444
$ws.Range("K94").Value = 13441.667
$ws.Range("L94").Value = 13908.444
$ws.Range("M94").Value = -12990.667
$ws.Range("N94").Value = -14810.444
$ws.Range("H105").Value = 10891.218
$ws.Range("I105").Value = 12174.315
$ws.Range("J105").Value = 4796.5
$ws.Range("K105").Value = 12174.315
$ws.Range("L105").Value = 4796.5
$ws.Range("M105").Value = -10427.315
$ws.Range("N105").Value = -8290.5
$ws.Range("H113").Value = 6068.5264
$ws.Range("I113").Value = 6884.5
$ws.Range("J113").Value = 1716.6666
$ws.Range("K113").Value = 6884.5
$ws.Range("L113").Value = 1716.6666
$ws.Range("M113").Value = -4714.5
$ws.Range("N113").Value = -6056.6666
$ws.Range("H134").Value = 15750.5
$ws.Range("I134").Value = 3998.3333
$ws.Range("J134").Value = 51007
$ws.Range("K134").Value = 11994.9999
$ws.Range("L134").Value = 153021
$ws.Range("M134").Value = -9459.999899999999
$ws.Range("N134").Value = -158091

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 6260.8
$ws.Range("I126").Value = 7319.5
$ws.Range("J126").Value = 4143.4
$ws.Range("K126").Value = 21958.5
$ws.Range("L126").Value = 12430.2
$ws.Range("M126").Value = -19488.5
$ws.Range("N126").Value = -17370.2

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 426774.44
$ws.Range("I136").Value = 1962.4667
$ws.Range("J136").Value = 1063992.4
$ws.Range("K136").Value = 5887.4001
$ws.Range("L136").Value = 3191977.2
$ws.Range("M136").Value = -3337.4001
$ws.Range("N136").Value = -3197077.2
